$wb = $excel.ActiveWorkbook

# --- Sheet 1: PERMANOVA ---
$ws1 = $wb.Worksheets.Item("PERMANOVA")

$ws1.Range("B2").Value = 0.3783812185883571
$ws1.Range("C2").Value = 0.05841353132863106
$ws1.Range("D2").Value = 2.731812415818474
$ws1.Range("E2").Value = 0.0043

$ws1.Range("A3").Value = 9
$ws1.Range("B3").Value = 3.052045430404583
$ws1.Range("C3").Value = 0.4711670204733284
$ws1.Range("D3").Value = 2.448329060401857
$ws1.Range("E3").Value = 0.0001

$ws1.Range("A4").Value = 22
$ws1.Range("B4").Value = 3.047202934118666
$ws1.Range("C4").Value = 0.4704194481980408

$ws1.Range("B5").Value = 6.477629583111605

# --- Sheet 2: PERMDISP ---
$ws2 = $wb.Worksheets.Item("PERMDISP")

$ws2.Range("B2").Value = 0.00345141700045872
$ws2.Range("C2").Value = 0.00345141700045872
$ws2.Range("D2").Value = 0.3673203149951823
$ws2.Range("F2").Value = 0.546

$ws2.Range("B3").Value = 0.2912823566962901
$ws2.Range("C3").Value = 0.009396205054719034
